$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I62").Value = 3000
$ws.Range("L62").Value = 5291.25
$ws.Range("K62").Value = 3000
$ws.Range("J62").Value = 5291.25
$ws.Range("H62").Value = 4309.2856
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -6539.25
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 26456.25
$ws.Range("H65").Value = 4309.2856
$ws.Range("I65").Value = 3000
$ws.Range("N65").Value = -32696.25
$ws.Range("J65").Value = 5291.25
$ws.Range("M65").Value = -11880
$ws.Range("L70").Value = 8805.125100000001
$ws.Range("J70").Value = 2935.0417
$ws.Range("N70").Value = -9345.125100000001
$ws.Range("H70").Value = 2756.3333
$ws.Range("J73").Value = 2935.0417
$ws.Range("H73").Value = 2756.3333
$ws.Range("L73").Value = 8805.125100000001
$ws.Range("N73").Value = -10677.1251
$ws.Range("M138").ClearContents()
$ws.Range("H138").Value = 7276.1924
$ws.Range("J138").Value = 7276.1924
$ws.Range("L138").Value = 21828.5772
$ws.Range("N138").Value = -32108.5772
$ws.Range("K138").Value = 0
$ws.Range("I138").Value = 0

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("I45").Value = 974.1429000000001
$ws.Range("K45").Value = 974.1429000000001
$ws.Range("H45").Value = 1484.4615
$ws.Range("M45").Value = -597.1429000000001
$ws.Range("K61").Value = 759.89655
$ws.Range("L61").Value = 3563.6553
$ws.Range("N61").Value = -3987.6553
$ws.Range("H61").Value = 2161.776
$ws.Range("J61").Value = 3563.6553
$ws.Range("M61").Value = -547.89655
$ws.Range("I61").Value = 759.89655
$ws.Range("H136").Value = 2161.776
$ws.Range("I136").Value = 759.89655
$ws.Range("L136").Value = 10690.9659
$ws.Range("J136").Value = 3563.6553
$ws.Range("K136").Value = 2279.68965
$ws.Range("N136").Value = -15790.9659
$ws.Range("M136").Value = 270.3103499999997

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M25").Value = -3064.6667
$ws.Range("H25").Value = 29986.2
$ws.Range("K25").Value = 3299.6667
$ws.Range("I25").Value = 3299.6667
$ws.Range("K64").Value = 470.66666
$ws.Range("M64").Value = -245.66666
$ws.Range("J64").Value = 540.46155
$ws.Range("N64").Value = -990.46155
$ws.Range("I64").Value = 470.66666
$ws.Range("H64").Value = 527.375
$ws.Range("L64").Value = 540.46155
$ws.Range("L67").Value = 540.46155
$ws.Range("I67").Value = 470.66666
$ws.Range("N67").Value = -2100.46155
$ws.Range("K67").Value = 470.66666
$ws.Range("H67").Value = 527.375
$ws.Range("J67").Value = 540.46155
$ws.Range("M67").Value = 309.33334
$ws.Range("N94").Value = -1764.5
$ws.Range("M94").Value = -123.0833
$ws.Range("K94").Value = 574.0833
$ws.Range("I94").Value = 574.0833
$ws.Range("H94").Value = 689.45
$ws.Range("J94").Value = 862.5
$ws.Range("L94").Value = 862.5
$ws.Range("J134").Value = 4628.5713
$ws.Range("N134").Value = -18955.7139
$ws.Range("H134").Value = 3320.543
$ws.Range("L134").Value = 13885.7139
$ws.Range("K134").Value = 8980.606800000001
$ws.Range("I134").Value = 2993.5356
$ws.Range("M134").Value = -6445.606800000001
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J16").Value = 5000
$ws.Range("M16").Value = -1085.5
$ws.Range("N16").Value = -5574
$ws.Range("H16").Value = 2927.1428
$ws.Range("I16").Value = 1372.5
$ws.Range("L16").Value = 5000
$ws.Range("K16").Value = 1372.5
$ws.Range("H31").Value = 2272.6812
$ws.Range("M31").Value = -1133.2307
$ws.Range("N31").Value = -3960.4666
$ws.Range("L31").Value = 3370.4666
$ws.Range("K31").Value = 1428.2307
$ws.Range("I31").Value = 1428.2307
$ws.Range("J31").Value = 3370.4666
$ws.Range("M34").Value = -1226.2307
$ws.Range("N34").Value = -3774.4666
$ws.Range("K34").Value = 1428.2307
$ws.Range("J34").Value = 3370.4666
$ws.Range("H34").Value = 2272.6812
$ws.Range("I34").Value = 1428.2307
$ws.Range("L34").Value = 3370.4666
$ws.Range("J113").Value = 5000
$ws.Range("N113").Value = -9340
$ws.Range("H113").Value = 2927.1428
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 797.5
$ws.Range("K113").Value = 1372.5
$ws.Range("I113").Value = 1372.5
$ws.Range("I132").Value = 2328.5
$ws.Range("K132").Value = 6985.5
$ws.Range("M132").Value = -4455.5
$ws.Range("H132").Value = 3946.1765
$ws.Range("H134").Value = 1864.2037
$ws.Range("K134").Value = 4334.634
$ws.Range("I134").Value = 1444.878
$ws.Range("M134").Value = -1799.634

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3175.2222
$ws.Range("L68").Value = 9525.6666
$ws.Range("K68").Value = 0
$ws.Range("J68").Value = 3175.2222
$ws.Range("I68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -11147.6666
$ws.Range("M71").ClearContents()
$ws.Range("L71").Value = 28576.9998
$ws.Range("K71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("N71").Value = -36688.99980000001
$ws.Range("J71").Value = 3175.2222
$ws.Range("H71").Value = 3175.2222
$ws.Range("H107").Value = 824.6197
$ws.Range("L107").Value = 3545.7
$ws.Range("M107").Value = 230.4146999999998
$ws.Range("J107").Value = 1181.9
$ws.Range("K107").Value = 1689.5853
$ws.Range("N107").Value = -7385.700000000001
$ws.Range("I107").Value = 563.1951

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M20").ClearContents()
$ws.Range("K20").Value = 0
$ws.Range("H20").Value = 70006
$ws.Range("I20").Value = 0
$ws.Range("H135").Value = 23794.285
$ws.Range("J135").Value = 23794.285
$ws.Range("L135").Value = 23794.285
$ws.Range("N135").Value = -33934.285

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N40").Value = -4938.6665
$ws.Range("M40").Value = -1316
$ws.Range("L40").Value = 4666.6665
$ws.Range("H40").Value = 3380.8
$ws.Range("J40").Value = 4666.6665
$ws.Range("I40").Value = 1452
$ws.Range("K40").Value = 1452
$ws.Range("N69").ClearContents()
$ws.Range("H69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("N76").ClearContents()
$ws.Range("J76").Value = 0
$ws.Range("H76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("J79").Value = 0
$ws.Range("H82").Value = 3781.2942
$ws.Range("I82").Value = 2468.7778
$ws.Range("N82").Value = -5979.875
$ws.Range("K82").Value = 2468.7778
$ws.Range("J82").Value = 5257.875
$ws.Range("L82").Value = 5257.875
$ws.Range("M82").Value = -2107.7778
$ws.Range("N85").Value = -7753.875
$ws.Range("I85").Value = 2468.7778
$ws.Range("K85").Value = 2468.7778
$ws.Range("J85").Value = 5257.875
$ws.Range("H85").Value = 3781.2942
$ws.Range("L85").Value = 5257.875
$ws.Range("M85").Value = -1220.7778
$ws.Range("H86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H88").Value = 29800
$ws.Range("L88").Value = 29800
$ws.Range("N88").Value = -30656
$ws.Range("J88").Value = 29800
$ws.Range("H89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("N91").Value = -32764
$ws.Range("L91").Value = 29800
$ws.Range("H91").Value = 29800
$ws.Range("J91").Value = 29800
$ws.Range("M93").Value = 396.5333000000001
$ws.Range("H93").Value = 1255.8096
$ws.Range("J93").Value = 2266.6667
$ws.Range("I93").Value = 851.4666999999999
$ws.Range("L93").Value = 2266.6667
$ws.Range("N93").Value = -4762.6667
$ws.Range("K93").Value = 851.4666999999999
$ws.Range("I132").Value = 1208.0435
$ws.Range("N132").Value = -14003.2499
$ws.Range("K132").Value = 3624.1305
$ws.Range("L132").Value = 8943.249899999999
$ws.Range("J132").Value = 2981.0833
$ws.Range("M132").Value = -1094.1305
$ws.Range("H132").Value = 1815.9429
$ws.Range("H135").Value = 30000
$ws.Range("J135").Value = 30000
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I132").Value = 3655.2727
$ws.Range("N132").Value = -46379.546
$ws.Range("K132").Value = 10965.8181
$ws.Range("L132").Value = 41319.546
$ws.Range("J132").Value = 13773.182
$ws.Range("M132").Value = -8435.8181
$ws.Range("H132").Value = 8714.227999999999
